$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2416933
$ws.Range("B3").Value = 100400
$ws.Range("B4").Value = 19.9
$ws.Range("B5").Value = 2.255
